# Interview_questions.xlsx edit script
# Commit message: "uploded JS Interview and DSA questions"
#
# Summary of changes:
#  1. Rename sheet "Coding" -> "DSA"
#  2. Javascript sheet: add note in B12 (deep-clone question) next to
#     "Shallow copy and deep copy of object", wrap text + taller row.
#  3. DSA sheet: add a Google-Sheets link in B1 next to the header cell.
#  4. Various selection / active-cell bookkeeping left behind by the
#     author while navigating the workbook (cosmetic, but reproduced for
#     fidelity): Typescript, express, State Management sheets gained a
#     "best fit" column (the user widened a column by double-clicking the
#     border after typing in column B), and the active sheet/tab ended up
#     on "nodejs".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New shared-string-backed content.
#    Order matters: the author typed the deep-clone note on the
#    Javascript sheet first, then pasted the Google Sheets link onto the
#    (renamed) DSA sheet, so write them in that order.
# ---------------------------------------------------------------------

$wsJs = $wb.Worksheets.Item("Javascript")
$wsJs.Range("B12").Value = "1. Ways to create a deep clone, does destructing creats a deep clones"
$wsJs.Range("B12").WrapText = $true
$wsJs.Rows.Item(12).RowHeight = 43.5

$wsCoding = $wb.Worksheets.Item("Coding")
$wsCoding.Range("B1").Value = "https://docs.google.com/spreadsheets/d/1hXserPuxVoWMG9Hs7y8wVdRCJTcj3xMBAEYUOXQ5Xag/edit?usp=drivesdk"

# ---------------------------------------------------------------------
# 2) Rename "Coding" -> "DSA"
# ---------------------------------------------------------------------
$wsCoding.Name = "DSA"
$wsDsa = $wb.Worksheets.Item("DSA")

# ---------------------------------------------------------------------
# 3) Column widths ("best fit" after typing the new content - values
#    chosen so the serialized width lands as close as possible to
#    Excel's own best-fit result for this content/font).
# ---------------------------------------------------------------------
$wsDsa.Columns.Item(1).ColumnWidth = 31
$wsDsa.Columns.Item(2).ColumnWidth = 96.83333333333333

$wsTs = $wb.Worksheets.Item("Typescript")
$wsTs.Columns.Item(1).ColumnWidth = 17.666666666666668

$wsExpress = $wb.Worksheets.Item("express")
$wsExpress.Columns.Item(1).ColumnWidth = 28.666666666666668

$wsState = $wb.Worksheets.Item("State Management")
$wsState.Columns.Item(1).ColumnWidth = 13.166666666666666

# ---------------------------------------------------------------------
# 4) Selections left on each sheet by the editing session.
# ---------------------------------------------------------------------
$wsDsa.Range("C6").Select() | Out-Null
$wsJs.Range("B12").Select() | Out-Null
$wsTs.Range("A3").Select() | Out-Null
$wsExpress.Range("A7").Select() | Out-Null
$wsState.Range("B7").Select() | Out-Null

$wsSql = $wb.Worksheets.Item("SQL")
$wsSql.Range("A13").Select() | Out-Null

# ---------------------------------------------------------------------
# 5) Final active sheet/tab: the author ended the session on "nodejs".
# ---------------------------------------------------------------------
$wsNode = $wb.Worksheets.Item("nodejs")
$wsNode.Activate() | Out-Null
$wsNode.Range("A15").Select() | Out-Null
